# Scheduled-runner update: refresh market-board derived columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# on the Halicarnassus_Profits sheets with newly-pulled price data.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 147.55556
$ws.Cells.Item(2, 10).Value = 472
$ws.Cells.Item(2, 12).Value = 472
$ws.Cells.Item(2, 14).Value = -698
$ws.Cells.Item(55, 8).Value = 2099.9
$ws.Cells.Item(55, 9).Value = 966.8333
$ws.Cells.Item(55, 11).Value = 966.8333
$ws.Cells.Item(55, 13).Value = -752.8333
$ws.Cells.Item(100, 8).Value = 3161
$ws.Cells.Item(100, 9).Value = 3109
$ws.Cells.Item(100, 10).Value = 3265
$ws.Cells.Item(100, 11).Value = 3109
$ws.Cells.Item(100, 12).Value = 3265
$ws.Cells.Item(100, 13).Value = -2568
$ws.Cells.Item(100, 14).Value = -4347
$ws.Cells.Item(107, 8).Value = 319.83334
$ws.Cells.Item(107, 10).Value = 279.5
$ws.Cells.Item(107, 12).Value = 279.5
$ws.Cells.Item(107, 14).Value = -4119.5
$ws.Cells.Item(111, 8).Value = 1061.5
$ws.Cells.Item(111, 10).Value = 1746
$ws.Cells.Item(111, 12).Value = 5238
$ws.Cells.Item(111, 14).Value = -11372
$ws.Cells.Item(116, 8).Value = 3797.25
$ws.Cells.Item(116, 10).Value = 4000
$ws.Cells.Item(116, 12).Value = 4000
$ws.Cells.Item(116, 14).Value = -10884
$ws.Cells.Item(138, 8).Value = 4988.3335
$ws.Cells.Item(138, 9).Value = 459
$ws.Cells.Item(138, 11).Value = 1377
$ws.Cells.Item(138, 13).Value = 3763

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2923.2727
$ws.Cells.Item(61, 9).Value = 2923.2727
$ws.Cells.Item(61, 11).Value = 2923.2727
$ws.Cells.Item(61, 13).Value = -2711.2727
$ws.Cells.Item(74, 8).Value = 2320.6
$ws.Cells.Item(74, 9).Value = 2083.5454
$ws.Cells.Item(74, 10).Value = 2972.5
$ws.Cells.Item(74, 11).Value = 2083.5454
$ws.Cells.Item(74, 12).Value = 2972.5
$ws.Cells.Item(74, 13).Value = -1209.5454
$ws.Cells.Item(74, 14).Value = -4720.5
$ws.Cells.Item(77, 8).Value = 2320.6
$ws.Cells.Item(77, 9).Value = 2083.5454
$ws.Cells.Item(77, 10).Value = 2972.5
$ws.Cells.Item(77, 11).Value = 10417.727
$ws.Cells.Item(77, 12).Value = 14862.5
$ws.Cells.Item(77, 13).Value = -6049.726999999999
$ws.Cells.Item(77, 14).Value = -23598.5
$ws.Cells.Item(97, 8).Value = 864.6875
$ws.Cells.Item(97, 9).Value = 862
$ws.Cells.Item(97, 11).Value = 862
$ws.Cells.Item(97, 13).Value = -366
$ws.Cells.Item(101, 8).Value = 62639
$ws.Cells.Item(101, 10).Value = 62639
$ws.Cells.Item(101, 12).Value = 62639
$ws.Cells.Item(101, 14).Value = -69129
$ws.Cells.Item(104, 8).Value = 24987.5
$ws.Cells.Item(104, 10).Value = 24987.5
$ws.Cells.Item(104, 12).Value = 24987.5
$ws.Cells.Item(104, 14).Value = -31975.5
$ws.Cells.Item(122, 8).Value = 1228
$ws.Cells.Item(122, 9).Value = 1228
$ws.Cells.Item(122, 11).Value = 3684
$ws.Cells.Item(122, 13).Value = -1234
$ws.Cells.Item(132, 8).Value = 2031
$ws.Cells.Item(132, 9).Value = 1195
$ws.Cells.Item(132, 11).Value = 3585
$ws.Cells.Item(132, 13).Value = -1055
$ws.Cells.Item(136, 8).Value = 2923.2727
$ws.Cells.Item(136, 9).Value = 2923.2727
$ws.Cells.Item(136, 11).Value = 8769.8181
$ws.Cells.Item(136, 13).Value = -6219.8181
$ws.Cells.Item(137, 8).Value = 72333.164
$ws.Cells.Item(137, 10).Value = 80000
$ws.Cells.Item(137, 12).Value = 80000
$ws.Cells.Item(137, 14).Value = -90200

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(76, 8).Value = 10314
$ws.Cells.Item(76, 10).Value = 10314
$ws.Cells.Item(76, 12).Value = 10314
$ws.Cells.Item(76, 14).Value = -10944
$ws.Cells.Item(79, 8).Value = 10314
$ws.Cells.Item(79, 10).Value = 10314
$ws.Cells.Item(79, 12).Value = 10314
$ws.Cells.Item(79, 14).Value = -12498
$ws.Cells.Item(86, 8).Value = 6583.3335
$ws.Cells.Item(86, 9).Value = 3862.5
$ws.Cells.Item(86, 10).Value = 8760
$ws.Cells.Item(86, 11).Value = 3862.5
$ws.Cells.Item(86, 12).Value = 8760
$ws.Cells.Item(86, 13).Value = -2739.5
$ws.Cells.Item(86, 14).Value = -11006
$ws.Cells.Item(89, 8).Value = 6583.3335
$ws.Cells.Item(89, 9).Value = 3862.5
$ws.Cells.Item(89, 10).Value = 8760
$ws.Cells.Item(89, 11).Value = 19312.5
$ws.Cells.Item(89, 12).Value = 43800
$ws.Cells.Item(89, 13).Value = -13696.5
$ws.Cells.Item(89, 14).Value = -55032
$ws.Cells.Item(94, 8).Value = 500
$ws.Cells.Item(94, 9).Value = 500
$ws.Cells.Item(94, 11).Value = 500
$ws.Cells.Item(94, 13).Value = -49
$ws.Cells.Item(105, 8).Value = 1252608.4
$ws.Cells.Item(105, 9).Value = 2225303.8
$ws.Cells.Item(105, 10).Value = 2000
$ws.Cells.Item(105, 11).Value = 2225303.8
$ws.Cells.Item(105, 12).Value = 2000
$ws.Cells.Item(105, 13).Value = -2223556.8
$ws.Cells.Item(105, 14).Value = -5494

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2605.7693
$ws.Cells.Item(58, 9).Value = 931
$ws.Cells.Item(58, 11).Value = 931
$ws.Cells.Item(58, 13).Value = -728
$ws.Cells.Item(105, 8).Value = 2713.5
$ws.Cells.Item(105, 9).Value = 1539.4
$ws.Cells.Item(105, 10).Value = 4670.3335
$ws.Cells.Item(105, 11).Value = 1539.4
$ws.Cells.Item(105, 12).Value = 4670.3335
$ws.Cells.Item(105, 13).Value = 207.5999999999999
$ws.Cells.Item(105, 14).Value = -8164.3335
$ws.Cells.Item(132, 8).Value = 3600
$ws.Cells.Item(132, 9).Value = 3600
$ws.Cells.Item(132, 11).Value = 10800
$ws.Cells.Item(132, 13).Value = -8270
$ws.Cells.Item(136, 8).Value = 2605.7693
$ws.Cells.Item(136, 9).Value = 931
$ws.Cells.Item(136, 11).Value = 2793
$ws.Cells.Item(136, 13).Value = -243
$ws.Cells.Item(140, 8).Value = 108593.336
$ws.Cells.Item(140, 10).Value = 108593.336
$ws.Cells.Item(140, 12).Value = 108593.336
$ws.Cells.Item(140, 14).Value = -118953.336

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 8).Value = 133.66667
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(129, 8).Value = 1783.1818
$ws.Cells.Item(129, 10).Value = 1851.875
$ws.Cells.Item(129, 12).Value = 5555.625
$ws.Cells.Item(129, 14).Value = -15555.625
$ws.Cells.Item(132, 8).Value = 5751.3335
$ws.Cells.Item(132, 9).Value = 4749
$ws.Cells.Item(132, 10).Value = 6252.5
$ws.Cells.Item(132, 11).Value = 42741
$ws.Cells.Item(132, 12).Value = 56272.5
$ws.Cells.Item(132, 13).Value = -40211
$ws.Cells.Item(132, 14).Value = -61332.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 664
$ws.Cells.Item(97, 9).Value = 590.8461
$ws.Cells.Item(97, 10).Value = 822.5
$ws.Cells.Item(97, 11).Value = 590.8461
$ws.Cells.Item(97, 12).Value = 822.5
$ws.Cells.Item(97, 13).Value = -94.84609999999998
$ws.Cells.Item(97, 14).Value = -1814.5
$ws.Cells.Item(122, 8).Value = 2427.7
$ws.Cells.Item(122, 9).Value = 1853.3334
$ws.Cells.Item(122, 11).Value = 5560.0002
$ws.Cells.Item(122, 13).Value = -3110.0002

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2829.8333
$ws.Cells.Item(16, 9).Value = 2829.8333
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 2829.8333
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 14).Value = -2659.8333
$ws.Cells.Item(22, 8).Value = 1280
$ws.Cells.Item(22, 9).Value = 1683
$ws.Cells.Item(22, 10).Value = 1107.2858
$ws.Cells.Item(22, 11).Value = 1683
$ws.Cells.Item(22, 12).Value = 1107.2858
$ws.Cells.Item(22, 13).Value = -1388
$ws.Cells.Item(22, 14).Value = -1697.2858
$ws.Cells.Item(27, 8).Value = 1280
$ws.Cells.Item(27, 9).Value = 1683
$ws.Cells.Item(27, 10).Value = 1107.2858
$ws.Cells.Item(27, 11).Value = 1683
$ws.Cells.Item(27, 12).Value = 1107.2858
$ws.Cells.Item(27, 13).Value = -1576
$ws.Cells.Item(27, 14).Value = -1321.2858
$ws.Cells.Item(40, 8).Value = 6513.3335
$ws.Cells.Item(40, 9).Value = 6088.5713
$ws.Cells.Item(40, 10).Value = 8000
$ws.Cells.Item(40, 11).Value = 6088.5713
$ws.Cells.Item(40, 12).Value = 8000
$ws.Cells.Item(40, 13).Value = -5952.5713
$ws.Cells.Item(40, 14).Value = -8272
$ws.Cells.Item(68, 8).Value = 6772.5454
$ws.Cells.Item(68, 9).Value = 2999.6667
$ws.Cells.Item(68, 11).Value = 2999.6667
$ws.Cells.Item(68, 13).Value = -2250.6667
$ws.Cells.Item(71, 8).Value = 6772.5454
$ws.Cells.Item(71, 9).Value = 2999.6667
$ws.Cells.Item(71, 11).Value = 14998.3335
$ws.Cells.Item(71, 13).Value = -11254.3335
$ws.Cells.Item(93, 8).Value = 1389.8
$ws.Cells.Item(93, 9).Value = 1389.8
$ws.Cells.Item(93, 11).Value = 1389.8
$ws.Cells.Item(93, 13).Value = -141.8
$ws.Cells.Item(100, 8).Value = 7113.857
$ws.Cells.Item(100, 9).Value = 4949.5
$ws.Cells.Item(100, 11).Value = 4949.5
$ws.Cells.Item(100, 13).Value = -4408.5
$ws.Cells.Item(104, 8).Value = 23791.5
$ws.Cells.Item(104, 10).Value = 23791.5
$ws.Cells.Item(104, 12).Value = 23791.5
$ws.Cells.Item(104, 14).Value = -30779.5
$ws.Cells.Item(132, 8).Value = 3223.1538
$ws.Cells.Item(132, 9).Value = 3040.1
$ws.Cells.Item(132, 11).Value = 9120.299999999999
$ws.Cells.Item(132, 13).Value = -6590.299999999999
$ws.Cells.Item(136, 8).Value = 1895.7142
$ws.Cells.Item(136, 9).Value = 1895.7142
$ws.Cells.Item(136, 11).Value = 5687.142599999999
$ws.Cells.Item(136, 13).Value = -3137.142599999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 1499.3334
$ws.Cells.Item(96, 9).Value = 1499.3334
$ws.Cells.Item(96, 11).Value = 1499.3334
$ws.Cells.Item(96, 13).Value = -126.3334
$ws.Cells.Item(113, 8).Value = 1239.6364
$ws.Cells.Item(113, 9).Value = 1227.6
$ws.Cells.Item(113, 11).Value = 3682.8
$ws.Cells.Item(113, 13).Value = -1512.8
$ws.Cells.Item(122, 8).Value = 4114.353
$ws.Cells.Item(122, 9).Value = 3259.4546
$ws.Cells.Item(122, 10).Value = 5681.6665
$ws.Cells.Item(122, 11).Value = 9778.363799999999
$ws.Cells.Item(122, 12).Value = 17044.9995
$ws.Cells.Item(122, 13).Value = -7328.363799999999
$ws.Cells.Item(122, 14).Value = -21944.9995
$ws.Cells.Item(126, 8).Value = 5008.25
$ws.Cells.Item(126, 9).Value = 2119.8
$ws.Cells.Item(126, 10).Value = 7071.4287
$ws.Cells.Item(126, 11).Value = 6359.400000000001
$ws.Cells.Item(126, 12).Value = 21214.2861
$ws.Cells.Item(126, 13).Value = -3889.400000000001
$ws.Cells.Item(126, 14).Value = -26154.2861

# CUL row 98 (The Sweet Kiss of Death / Rice Vinegar): HQ leve no longer
# yields a distinct profit figure this cycle, so LeveProfitHQ is cleared.
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 14).ClearContents()
